$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3,1).Value = '"https://www.immowelt.de/expose/2aqpt5d"'
$ws.Cells.Item(3,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(3,3).Value = '"Die Einbauk\u00fcche ist schon drin - Sie m\u00fcssen nur noch einziehen! Gehobene 2-Z. Wohnung im Steinbachtal "'
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '1165'
$ws.Cells.Item(3,5).Value = '"EUR"'
$ws.Cells.Item(3,6).NumberFormat = "@"
$ws.Cells.Item(3,6).Value = '77'
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = '2'
$ws.Cells.Item(3,8).Value = 'null'
$ws.Cells.Item(3,9).Value = '["renoviert", "GAS", "Zentralheizung", "gartennutzung", "Einbauk\u00fcche", "Terrasse"]'
$ws.Cells.Item(3,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Steinbachtal", "ZipCode": "97082", "LocationId": 496023, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(3,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "3495", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Steinbachtal)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-steinbachtal/mietspiegel"}}, "DataTable": [{"NumberValue": 1165, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 280, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1445, "Unit": "EUR"}, {"NumberValue": 75, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "1 Stellplatz"}]}'
$ws.Cells.Item(3,12).Value = '["https://ms.immowelt.org/d789ca7c-5451-48cf-a649-7af7d4c9078a/0af4290a-753d-4df9-bd0e-3989d2db4215", "https://ms.immowelt.org/2908c42c-0f96-410b-b619-dba38e6a42a1/75b141df-d3a1-4754-b265-67306b3a36e1", "https://ms.immowelt.org/fbb8e8e1-0cc5-4a95-8902-c6f6af96acc8/548e065f-569a-4838-9887-a19794c56f49", "https://ms.immowelt.org/4365c939-93fc-4ec1-b3e8-77d586d0e200/5112a29b-5171-4dc0-b36a-b9b4c4ac411f", "https://ms.immowelt.org/db7ac2e7-2969-42be-920b-36be7e6aab0d/0cf7253b-4260-4efa-af51-0074b8782be6", "https://ms.immowelt.org/d1a4bcbb-b277-4731-8486-0a09fb1b8996/475444d0-b794-44bc-ba59-ce3f2f51b9f8", "https://ms.immowelt.org/748a1662-4a4a-4f6d-bb72-baf731de8adc/700fbb9a-6ef9-48a1-8e73-4c46cba5c258", "https://ms.immowelt.org/ed9af812-426e-4e49-8e89-2ede0c407269/5a107fdf-e72a-4f22-9de2-52351751a755", "https://ms.immowelt.org/643489b1-7bf2-46a3-882c-7d8b8b0957ca/bce02d3c-331a-4110-8e99-d872140474b1", "https://ms.immowelt.org/d8490c94-8bc2-4f38-affe-7986af339d1c/5d00296c-3134-4d05-a29a-760985379123", "https://ms.immowelt.org/862887dd-fca5-4504-a254-d22277945d21/1fdaef5f-474b-4b44-84df-c23fb4960205", "https://ms.immowelt.org/235cb2e7-3834-406d-a4ac-0a1db7702757/b3d2810f-2282-453d-99bf-f70073155c71", "https://ms.immowelt.org/3ae43d7c-7b0d-435e-84bb-7a6ea4aafc03/ae2064f4-8fe0-463b-958f-a7c34e88466a"]'
$ws.Cells.Item(3,13).Value = '""'
$ws.Cells.Item(3,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Die Einbauk\u00fcche ist schon drin - Sie m\u00fcssen nur noch einziehen! Gehobene 2-Z. Wohnung im Steinbachtal ", "basicObjectPricEur": 1165, "basicLivingSpace": 77, "basicRooms": 2, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2aqpt5d", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 4
$ws.Cells.Item(4,1).Value = '"https://www.immowelt.de/expose/2tfuz4l"'
$ws.Cells.Item(4,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(4,3).Value = '"Bestlage mit Festungsblick"'
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '610'
$ws.Cells.Item(4,5).Value = '"EUR"'
$ws.Cells.Item(4,6).NumberFormat = "@"
$ws.Cells.Item(4,6).Value = '65'
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = '1'
$ws.Cells.Item(4,8).Value = 'null'
$ws.Cells.Item(4,9).Value = '["DUSCHE", "FENSTER", "Dachgeschoss", "frei", "offene K\u00fcche", "Stellplatz"]'
$ws.Cells.Item(4,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Sanderau", "ZipCode": "97072", "Street": "Breslauer Stra\u00dfe", "LocationId": 496022, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(4,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "Abl\u00f6se der K\u00fcche", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Sanderau)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg/mietspiegel"}}, "DataTable": [{"NumberValue": 610, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete", "Comments": ["zzgl. Nebenkosten"]}]}'
$ws.Cells.Item(4,12).Value = '["https://ms.immowelt.org/b4607fa6-96c6-40a3-8049-ed2221dafd3c/bb9b114c-9b74-4749-b16f-f5402ce89d13", "https://ms.immowelt.org/0678cd9c-641a-4bcc-a0ac-02f46fc8e1ea/b0aa6847-a6e4-474c-b300-23bb59c7628b", "https://ms.immowelt.org/d524e2d3-1330-4559-987d-6495b2a4c71d/3adebae2-58ee-44b1-a898-86e63f8eff12", "https://ms.immowelt.org/61eb287f-36c7-4a51-a4db-f88e4049ad30/29721137-7a70-4a6d-9abc-340c019663bf"]'
$ws.Cells.Item(4,13).Value = '""'
$ws.Cells.Item(4,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Bestlage mit Festungsblick", "basicObjectPricEur": 610, "basicLivingSpace": 65, "basicRooms": 1, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicStreet": "Breslauer Stra\u00dfe", "basicUrl": "https://www.immowelt.de/expose/2tfuz4l", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 5
$ws.Cells.Item(5,1).Value = '"https://www.immowelt.de/expose/2a8j75n"'
$ws.Cells.Item(5,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(5,3).Value = '"Wohnberechtigungsschein erforderlich - 2 Zi, 58m\u00b2, Balkon"'
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '554.99'
$ws.Cells.Item(5,5).Value = '"EUR"'
$ws.Cells.Item(5,6).NumberFormat = "@"
$ws.Cells.Item(5,6).Value = '58.42'
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = '2'
$ws.Cells.Item(5,8).NumberFormat = "@"
$ws.Cells.Item(5,8).Value = '2019'
$ws.Cells.Item(5,9).Value = '["Personenaufzug", "Stellplatz"]'
$ws.Cells.Item(5,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 15", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(5,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1660", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 554.99, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 210, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 208, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1032.99, "Unit": "EUR"}]}'
$ws.Cells.Item(5,12).Value = '["https://ms.immowelt.org/2c287869-b002-4a19-baea-0433bdc7a327/e15bc071-fe33-4461-8b34-9efe6f7f6e23", "https://ms.immowelt.org/5b165c7f-f365-48c1-841e-a8b1fd2bb6a1/23a21cb8-5f33-4e46-b102-f4960458c3e9", "https://ms.immowelt.org/4b517be9-8d93-492f-ac56-516e6bad3d95/1238a7b6-1272-4e53-9ead-1eafab0548c7", "https://ms.immowelt.org/acdfaf67-7ef9-4ee5-a2e3-416020180ea6/b120a8ca-7b3e-47f0-99b1-006c6eb00277", "https://ms.immowelt.org/540c5813-745e-419c-9217-36ba5261e842/259d75f3-7aca-466d-b4ed-f4b5a08ff770", "https://ms.immowelt.org/9ac46b17-c8c4-446e-8e96-de026ce267b4/f6098bb6-68df-4742-8eda-56e4660c3221", "https://ms.immowelt.org/09537047-f692-4086-b911-c08acb726d74/595eae7d-cd2e-4c30-b4d5-6bf58c421952"]'
$ws.Cells.Item(5,13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "Herr", "firstName": "Stefan", "lastName": "Brogl"}'
$ws.Cells.Item(5,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnberechtigungsschein erforderlich - 2 Zi, 58m\u00b2, Balkon", "basicObjectPricEur": 554.99, "basicLivingSpace": 58.42, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 15", "basicUrl": "https://www.immowelt.de/expose/2a8j75n"}'

# Row 6
$ws.Cells.Item(6,1).Value = '"https://www.immowelt.de/expose/2aaj75n"'
$ws.Cells.Item(6,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(6,3).Value = '"Wohnung kann nur mit Wohnberechtigungsschein angemietet werden! 2 Zimmer auf 56 m\u00b2"'
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '532.19'
$ws.Cells.Item(6,5).Value = '"EUR"'
$ws.Cells.Item(6,6).NumberFormat = "@"
$ws.Cells.Item(6,6).Value = '56.02'
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = '2'
$ws.Cells.Item(6,8).NumberFormat = "@"
$ws.Cells.Item(6,8).Value = '2019'
$ws.Cells.Item(6,9).Value = '["Personenaufzug"]'
$ws.Cells.Item(6,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 17", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(6,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1590", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 532.19, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 179, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 110.92, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 882.11, "Unit": "EUR"}]}'
$ws.Cells.Item(6,12).Value = '["https://ms.immowelt.org/814cec20-d48a-4483-ae7f-5b779e29c323/404d9d3f-baa2-473d-af13-096d09772807", "https://ms.immowelt.org/60ba3310-b6df-446e-b860-6deb11ef0a09/2558c572-051e-464d-9aad-205683e0f215", "https://ms.immowelt.org/e10fa924-843a-4cd9-897d-8a60a9cd8473/3fd9082b-7a00-4210-99b9-f828a6a2e085", "https://ms.immowelt.org/20e41405-47ae-47ac-97d8-0e7340de1982/b112096a-8211-4248-a5e6-0904121ddac0", "https://ms.immowelt.org/2300c3d0-702d-4867-94fc-810ec8ca8971/f6adb669-5e8c-4182-ad4f-a5b6f477b025", "https://ms.immowelt.org/bdd38a88-61e2-4f96-bac0-895a2e862eed/9b1e1d00-b895-496c-bf47-2912364590ab", "https://ms.immowelt.org/a047e5df-c99a-4420-a028-e84ef22f3786/0dd57565-cb4f-4092-ab20-2a61b660568a"]'
$ws.Cells.Item(6,13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "Herr", "firstName": "Stefan", "lastName": "Brogl"}'
$ws.Cells.Item(6,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnung kann nur mit Wohnberechtigungsschein angemietet werden! 2 Zimmer auf 56 m\u00b2", "basicObjectPricEur": 532.19, "basicLivingSpace": 56.02, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 17", "basicUrl": "https://www.immowelt.de/expose/2aaj75n"}'

# Row 7
$ws.Cells.Item(7,1).Value = '"https://www.immowelt.de/expose/2a97e5s"'
$ws.Cells.Item(7,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(7,3).Value = '"Gut gemacht !! 3 Zi, 81 qm, Balkon und Einbauk\u00fcche"'
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '1183.37'
$ws.Cells.Item(7,5).Value = '"EUR"'
$ws.Cells.Item(7,6).NumberFormat = "@"
$ws.Cells.Item(7,6).Value = '81.32'
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = '3'
$ws.Cells.Item(7,8).NumberFormat = "@"
$ws.Cells.Item(7,8).Value = '2019'
$ws.Cells.Item(7,9).Value = '["Personenaufzug", "Einbauk\u00fcche", "Stellplatz"]'
$ws.Cells.Item(7,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 13", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(7,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "3550", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 1183.37, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 252, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 171, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1666.37, "Unit": "EUR"}]}'
$ws.Cells.Item(7,12).Value = '["https://ms.immowelt.org/4e4fecdb-3fe4-4b73-be6f-d02dcff6e907/a23bae17-4a7c-4155-88b0-8014e5e3021a", "https://ms.immowelt.org/324f63aa-2fff-4fc8-8ba1-ac69ea7a7e23/ed4f84b2-861b-4e99-b1c6-4e3963ae1f35", "https://ms.immowelt.org/be4e9cc8-6ad5-4ae2-a739-1b8b7178ee12/ba6c9877-b521-46d0-a222-9e1e70f8ffd4"]'
$ws.Cells.Item(7,13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "herr", "firstName": "Stefan", "lastName": "Brogl", "phone": "+49 931 306 990 20"}'
$ws.Cells.Item(7,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Gut gemacht !! 3 Zi, 81 qm, Balkon und Einbauk\u00fcche", "basicObjectPricEur": 1183.37, "basicLivingSpace": 81.32, "basicRooms": 3, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 13", "basicUrl": "https://www.immowelt.de/expose/2a97e5s", "basicContactPhone": "+49 931 306 990 20"}'

# Row 8
$ws.Cells.Item(8,1).Value = '"https://www.immowelt.de/expose/25dtq5g"'
$ws.Cells.Item(8,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(8,3).Value = '"3-ZW W\u00fcrzburg  Sanderau,  3. OG.  Bad neu."'
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '950'
$ws.Cells.Item(8,5).Value = '"EUR"'
$ws.Cells.Item(8,6).NumberFormat = "@"
$ws.Cells.Item(8,6).Value = '79'
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = '3'
$ws.Cells.Item(8,8).NumberFormat = "@"
$ws.Cells.Item(8,8).Value = '1911'
$ws.Cells.Item(8,9).Value = '["Altbau (bis 1945)", "WANNE", "FENSTER", "Kelleranteil", "GAS", "Etagenheizung", "Einbauk\u00fcche", "Balkon"]'
$ws.Cells.Item(8,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Sanderau", "ZipCode": "97072", "Street": "Eichendorffstrasse 4", "LocationId": 496022, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(8,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1900", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Sanderau)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg/mietspiegel"}}, "DataTable": [{"NumberValue": 950, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 120, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}]}'
$ws.Cells.Item(8,12).Value = '["https://ms.immowelt.org/53a40e8a-24fb-45ae-a531-fff2a2f4ef38/73c3faa7-d432-44a3-995c-dd6a57b55ff8", "https://ms.immowelt.org/5fd3d699-98a8-4599-a121-418d851bdca3/552eef5e-2a9f-4c01-b006-8b517f701ec2", "https://ms.immowelt.org/e1c22ef7-fa4f-4bbe-b945-4a0c9417515d/e28a4ad5-ff43-4cdd-a230-41a8a963e3de", "https://ms.immowelt.org/1fe9ddd3-4952-4f94-9aba-7b98e8a565e2/4cc8d480-3a53-414c-ada9-ef0da80900ed", "https://ms.immowelt.org/2e8ae348-cab2-4e79-8266-ee52099c1e73/524bcb38-9f3e-44ef-afa6-8810bc8ec8bf", "https://ms.immowelt.org/b0e64c11-6bbf-4ed6-9eb5-3767349429ee/e3628c5a-5b53-4019-a135-ff500cf7a0ac", "https://ms.immowelt.org/4a8c3d49-12ed-417a-b206-b38a38d9e894/7e1402ca-1074-4126-8ddc-ba9dbe8eac83", "https://ms.immowelt.org/82deeace-1bd0-4723-8a33-2c2ad6084c4a/cabe6e07-e983-4c08-826a-e189aa3d4288", "https://ms.immowelt.org/4e16951a-d91d-485c-a019-726374d3e5d3/97f13c47-29c0-4a36-8ee3-73268762ea32", "https://ms.immowelt.org/f8cc63cd-6dc9-4cc3-8aba-e6dbc1f8c911/bebff58a-a1df-4e6a-868b-ba584b296f85"]'
$ws.Cells.Item(8,13).Value = '""'
$ws.Cells.Item(8,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "3-ZW W\u00fcrzburg  Sanderau,  3. OG.  Bad neu.", "basicObjectPricEur": 950, "basicLivingSpace": 79, "basicRooms": 3, "basicConstructionYear": 1911, "basicCity": "W\u00fcrzburg", "basicStreet": "Eichendorffstrasse 4", "basicUrl": "https://www.immowelt.de/expose/25dtq5g", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 9
$ws.Cells.Item(9,1).Value = '"https://www.immowelt.de/expose/2ajt75p"'
$ws.Cells.Item(9,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(9,3).Value = '"Wohnberechtigungsschein (EOF3) erforderlich!! 2-Zimmer Wohnung mit Dachterrasse!!"'
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '474.53'
$ws.Cells.Item(9,5).Value = '"EUR"'
$ws.Cells.Item(9,6).NumberFormat = "@"
$ws.Cells.Item(9,6).Value = '49.95'
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = '2'
$ws.Cells.Item(9,8).NumberFormat = "@"
$ws.Cells.Item(9,8).Value = '2019'
$ws.Cells.Item(9,9).Value = '["Personenaufzug"]'
$ws.Cells.Item(9,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 15", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(9,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1420", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 474.53, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 167, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 98.9, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 740.43, "Unit": "EUR"}]}'
$ws.Cells.Item(9,12).Value = '["https://ms.immowelt.org/059062b5-c451-4105-bfa0-dfb0db8ca29f/0f1cbcbd-e9ed-45e6-96e8-910e278b3044", "https://ms.immowelt.org/ca46e538-90d4-4fa2-8f22-2b6881ac1800/9a457fdc-44b3-4336-b20e-89b86824b294", "https://ms.immowelt.org/a1a88b02-0504-429f-b458-7c1a33126061/b8bc3487-e427-4a1b-b30b-2df51d481f90", "https://ms.immowelt.org/d507b7b4-5e3c-452e-957a-b5636948e585/d2a53e80-0b41-44e4-a681-ce8531a32c32", "https://ms.immowelt.org/ff43f21c-a328-4711-9ce2-41a9bff1bf59/83bd2e55-f9e8-44b1-bf30-e874b0a81423"]'
$ws.Cells.Item(9,13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "Herr", "firstName": "Stefan", "lastName": "Brogl"}'
$ws.Cells.Item(9,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnberechtigungsschein (EOF3) erforderlich!! 2-Zimmer Wohnung mit Dachterrasse!!", "basicObjectPricEur": 474.53, "basicLivingSpace": 49.95, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 15", "basicUrl": "https://www.immowelt.de/expose/2ajt75p"}'

# Row 10
$ws.Cells.Item(10,1).Value = '"https://www.immowelt.de/expose/2a7j75n"'
$ws.Cells.Item(10,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(10,3).Value = '"Wohnen im herrlichen HUBland: 2Zi, 70m\u00b2, EBK und Balkon"'
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '1181.21'
$ws.Cells.Item(10,5).Value = '"EUR"'
$ws.Cells.Item(10,6).NumberFormat = "@"
$ws.Cells.Item(10,6).Value = '70.87'
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = '2'
$ws.Cells.Item(10,8).NumberFormat = "@"
$ws.Cells.Item(10,8).Value = '2019'
$ws.Cells.Item(10,9).Value = '["Personenaufzug", "Einbauk\u00fcche", "Stellplatz"]'
$ws.Cells.Item(10,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 13", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(10,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "3540", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 1181.21, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 219, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 169, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1629.21, "Unit": "EUR"}]}'
$ws.Cells.Item(10,12).Value = '["https://ms.immowelt.org/0d9ec1ef-21cc-484a-88ae-193c7267a3bb/1e78105c-6490-4468-a6c3-936b8b388b80", "https://ms.immowelt.org/5c55f670-770f-44e9-b26b-7c9164f4a76a/9f56e225-8915-45d3-9af7-0e36ad2d804e", "https://ms.immowelt.org/12a9fd2a-da00-49a5-9999-efb246a9d727/efc35e4b-a860-47d1-823d-1282d60e0bc6", "https://ms.immowelt.org/bb9a6b97-0d40-46d1-9b24-011a5dcc7508/40526838-218d-4a4c-9422-8bd98023d9f6", "https://ms.immowelt.org/1a701923-58c4-4074-92b3-3e6bc100d51c/2f4849ee-2b95-432d-8ed4-b41b0211a117", "https://ms.immowelt.org/4d6aee50-6e8f-423f-b67a-bd01a011592d/4d58cead-e166-4a77-8c0f-2ee869a9e824"]'
$ws.Cells.Item(10,13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "herr", "firstName": "Stefan", "lastName": "Brogl", "phone": "+49 931 306 990 20"}'
$ws.Cells.Item(10,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnen im herrlichen HUBland: 2Zi, 70m\u00b2, EBK und Balkon", "basicObjectPricEur": 1181.21, "basicLivingSpace": 70.87, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 13", "basicUrl": "https://www.immowelt.de/expose/2a7j75n", "basicContactPhone": "+49 931 306 990 20"}'

# Row 11
$ws.Cells.Item(11,1).Value = '"https://www.immowelt.de/expose/2an6q5g"'
$ws.Cells.Item(11,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(11,3).Value = '"Wohnberechtigungsschein (WBS3) erforderlich!! Gem\u00fctliche Zweiraumwohnung auf 59m\u00b2"'
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '566.96'
$ws.Cells.Item(11,5).Value = '"EUR"'
$ws.Cells.Item(11,6).NumberFormat = "@"
$ws.Cells.Item(11,6).Value = '59.68'
$ws.Cells.Item(11,7).NumberFormat = "@"
$ws.Cells.Item(11,7).Value = '2'
$ws.Cells.Item(11,8).NumberFormat = "@"
$ws.Cells.Item(11,8).Value = '2019'
$ws.Cells.Item(11,9).Value = '["Stellplatz"]'
$ws.Cells.Item(11,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 15", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(11,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1700", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 566.96, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 199, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 124, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 949.96, "Unit": "EUR"}]}'
$ws.Cells.Item(11,12).Value = '["https://ms.immowelt.org/04bb8e8f-62a6-4aa1-be23-fcac75a2ca80/64e8eae5-448b-4e0a-aef6-3ee4487d6cad", "https://ms.immowelt.org/07409987-f4b5-4a2f-9ed5-55d000afa646/e6d2fc11-39c9-44c4-9ab7-6b79a5e4667d", "https://ms.immowelt.org/d8d286af-d70c-4388-a156-0e4525e576c2/3f08e4d7-3ff5-4332-b507-dfc163a7d33f", "https://ms.immowelt.org/0eb043d6-e29a-4797-8d47-b95fdf183238/bda6c2c7-7309-497f-836d-6483f70abd43", "https://ms.immowelt.org/4aed45df-9c7d-4a2b-b188-b8968b22727d/50fd789e-2343-4b28-867b-f726db140de4"]'
$ws.Cells.Item(11,13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "Herr", "firstName": "Stefan", "lastName": "Brogl"}'
$ws.Cells.Item(11,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnberechtigungsschein (WBS3) erforderlich!! Gem\u00fctliche Zweiraumwohnung auf 59m\u00b2", "basicObjectPricEur": 566.96, "basicLivingSpace": 59.68, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 15", "basicUrl": "https://www.immowelt.de/expose/2an6q5g"}'

# Row 12
$ws.Cells.Item(12,1).Value = '"https://www.immowelt.de/expose/2awv45u"'
$ws.Cells.Item(12,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(12,3).Value = '"Luxuri\u00f6se Penthouse-Wohnung Erstbezug! 3-Zimmer-Neubau in W\u00fcrzburg-Heidingfeld"'
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '1720'
$ws.Cells.Item(12,5).Value = '"EUR"'
$ws.Cells.Item(12,6).NumberFormat = "@"
$ws.Cells.Item(12,6).Value = '123.36'
$ws.Cells.Item(12,7).NumberFormat = "@"
$ws.Cells.Item(12,7).Value = '3'
$ws.Cells.Item(12,8).NumberFormat = "@"
$ws.Cells.Item(12,8).Value = '2022'
$ws.Cells.Item(12,9).Value = '["Neubau", "barriefrei", "WANNE", "gaestewc", "FENSTER", "Kelleranteil", "Personenaufzug", "FERNE", "GAS", "Fu\u00dfbodenheizung", "Zentralheizung", "rollstuhlgerecht", "PARKETT", "frei", "Einbauk\u00fcche", "offene K\u00fcche", "Speisekammer", "Stellplatz", "Tiefgarage", "REINIGUNG", "Balkon", "Terrasse"]'
$ws.Cells.Item(12,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Heidingsfeld", "ZipCode": "97084", "LocationId": 496013, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(12,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "5160", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Heidingsfeld)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-heidingsfeld/mietspiegel"}}, "DataTable": [{"NumberValue": 1720, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 320, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 2040, "Unit": "EUR"}, {"NumberValue": 70, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "2 Stellpl\u00e4tze, je"}]}'
$ws.Cells.Item(12,12).Value = '["https://ms.immowelt.org/70edfd50-f060-4619-8d8e-0a42b4111160/271e31b7-8484-45a8-8ac9-904dabf13c7b", "https://ms.immowelt.org/efa2a84b-bb24-4421-ae15-0649647fc5f0/8490cb37-3682-446f-8b11-e9e1da37cba9", "https://ms.immowelt.org/1683eefc-9ca0-4cbc-992f-8e3a9b044b29/7df67d66-b166-4c09-b89c-52d2409839a6", "https://ms.immowelt.org/80c22a83-a553-4a11-b6d3-845ab5aee5fa/a5ed848d-87ec-43a7-91b1-3a82889e2f11", "https://ms.immowelt.org/912ad1a7-5c04-49c1-bfdb-52dfd32f4078/087254f6-99be-4f33-ad19-a0d31350b6df", "https://ms.immowelt.org/c04aac07-656d-484f-97c8-f9183167c740/27c413f7-27b3-4874-9d22-7b3ce7fb1772", "https://ms.immowelt.org/6355ac0c-0919-41d1-9939-c7327ad98472/b5a33c7f-e707-4bd6-839a-ea939f810d7c", "https://ms.immowelt.org/e86eec55-1463-4303-8c06-6e1cd90ac14e/01d5e414-8279-49e7-ac42-160f28ed6a9b", "https://ms.immowelt.org/f6c4b53e-148d-44ad-b8fb-50e2e7c095c6/9403fef4-e110-4d99-88ec-08c732282eda", "https://ms.immowelt.org/bb63f6b4-e5b8-4f84-b1e9-c19308a3d41f/20e645fd-f306-449c-82ae-d3591730cca2", "https://ms.immowelt.org/6eff25e1-ca40-4fdc-b6f2-adcf3b4d0a45/63b7a19b-f618-4933-adb0-f0cd5e1263a5", "https://ms.immowelt.org/2811d6e1-5985-40bf-a07a-594fdaf52254/4a7a045f-966f-4645-86ce-900f6e22d446", "https://ms.immowelt.org/aefa4aed-f85d-4406-ae1b-9c233616a69e/fd830bae-b6de-4525-8e19-1c72c8249092", "https://ms.immowelt.org/c0855678-63ef-4feb-a503-a8d0adf77171/3851bcbc-db22-481c-b205-d0b9fb09017f", "https://ms.immowelt.org/df7e819e-aa27-4bbf-b743-a9e2712de66b/0737cd62-424e-4825-afba-0a294eb938ec", "https://ms.immowelt.org/f47e219f-f548-45dd-a07d-bb94c2b60e97/cfc0cf18-e610-4657-b7ed-3ca18984267a"]'
$ws.Cells.Item(12,13).Value = '{"companyName": "Gute Bude Immobilien", "address": {"city": "W\u00fcrzburg", "zipCode": "97082", "street": "Zeller Str. 3c"}, "salutation": "Frau", "firstName": "Nelly", "lastName": "Gronau"}'
$ws.Cells.Item(12,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Luxuri\u00f6se Penthouse-Wohnung Erstbezug! 3-Zimmer-Neubau in W\u00fcrzburg-Heidingfeld", "basicObjectPricEur": 1720, "basicLivingSpace": 123.36, "basicRooms": 3, "basicConstructionYear": 2022, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2awv45u"}'

# Row 14
$ws.Cells.Item(14,1).Value = '"https://www.immowelt.de/expose/2at6f5v"'
$ws.Cells.Item(14,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(14,3).Value = '"Wohnung vermieten in W\u00fcrzburg"'
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '420'
$ws.Cells.Item(14,5).Value = '"EUR"'
$ws.Cells.Item(14,6).Value = '""'
$ws.Cells.Item(14,7).Value = '""'
$ws.Cells.Item(14,8).Value = 'null'
$ws.Cells.Item(14,9).Value = '["GAS"]'
$ws.Cells.Item(14,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Gromb\u00fchl", "ZipCode": "97080", "LocationId": 496011, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(14,11).Value = '{"AdditionalInformation": {"MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Gromb\u00fchl)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-grombuehl/mietspiegel"}}, "DataTable": [{"NumberValue": 420, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete", "Comments": ["zzgl. Nebenkosten"]}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 420, "Unit": "EUR"}]}'
$ws.Cells.Item(14,12).Value = '[]'
$ws.Cells.Item(14,13).Value = '""'
$ws.Cells.Item(14,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnung vermieten in W\u00fcrzburg", "basicObjectPricEur": 420, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2at6f5v", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 15
$ws.Cells.Item(15,1).Value = '"https://www.immowelt.de/expose/2afyn5t"'
$ws.Cells.Item(15,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(15,3).Value = '"2,5 Zimmer Wohnung beim Japanischen Garten in W\u00fcrzburg"'
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '750'
$ws.Cells.Item(15,5).Value = '"EUR"'
$ws.Cells.Item(15,6).NumberFormat = "@"
$ws.Cells.Item(15,6).Value = '65'
$ws.Cells.Item(15,7).NumberFormat = "@"
$ws.Cells.Item(15,7).Value = '2.5'
$ws.Cells.Item(15,8).NumberFormat = "@"
$ws.Cells.Item(15,8).Value = '1960'
$ws.Cells.Item(15,9).Value = '["Altbau (bis 1945)", "WANNE", "Kelleranteil", "renoviert", "GAS", "Zentralheizung", "Einbauk\u00fcche"]'
$ws.Cells.Item(15,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Zellerau", "ZipCode": "97082", "LocationId": 496026, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(15,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "2000", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Zellerau)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-zellerau/mietspiegel"}}, "DataTable": [{"NumberValue": 750, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 200, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 950, "Unit": "EUR"}]}'
$ws.Cells.Item(15,12).Value = '["https://ms.immowelt.org/7392265b-7589-4e96-8719-b289bb8ace01/19ed71d6-701a-4bbe-929c-a35f0d2ec94b", "https://ms.immowelt.org/15baddee-aaf7-4f0a-b016-7f47a07c4e37/48316143-5143-447a-8c57-f0b7e929d6fe", "https://ms.immowelt.org/027c0258-c7df-4398-b7bc-fce818803786/c5879894-c70e-4d36-87c5-2f3d1391a860", "https://ms.immowelt.org/512f5ff2-4810-463d-be09-2b1645c65d7d/c20634e3-424e-4fa7-adc1-7bb96e05a731", "https://ms.immowelt.org/0bc625b6-c3c3-4417-81a8-a0c000046701/97c3d523-e7c6-4b79-af4f-d9928909527c", "https://ms.immowelt.org/76a5fc74-de02-4ae9-aec8-48e4168cea64/423affd6-d43f-46c2-9187-ca30a02614d6", "https://ms.immowelt.org/39f67b57-3390-44ea-ab71-3cbb1b166a73/03cd65f9-4b01-46f7-832d-c4ec60e69999", "https://ms.immowelt.org/4aa9aeb4-7027-43a0-85f3-887e51080c31/b7c1b3ae-576c-44d3-a98d-260591a7c97b"]'
$ws.Cells.Item(15,13).Value = '""'
$ws.Cells.Item(15,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "2,5 Zimmer Wohnung beim Japanischen Garten in W\u00fcrzburg", "basicObjectPricEur": 750, "basicLivingSpace": 65, "basicRooms": 2.5, "basicConstructionYear": 1960, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2afyn5t", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 16
$ws.Cells.Item(16,1).Value = '"https://www.immowelt.de/expose/2agbf5s"'
$ws.Cells.Item(16,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(16,3).Value = '"Traumwohnung in Top Lage mit toller Weitsicht!"'
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '2880'
$ws.Cells.Item(16,5).Value = '"EUR"'
$ws.Cells.Item(16,6).NumberFormat = "@"
$ws.Cells.Item(16,6).Value = '195'
$ws.Cells.Item(16,7).NumberFormat = "@"
$ws.Cells.Item(16,7).Value = '6'
$ws.Cells.Item(16,8).NumberFormat = "@"
$ws.Cells.Item(16,8).Value = '2023'
$ws.Cells.Item(16,9).Value = '["Neubau", "barriefrei", "WANNE", "gaestewc", "FENSTER", "Kelleranteil", "LUFTWP", "Garten", "gartennutzung", "rollstuhlgerecht", "Balkon", "Terrasse"]'
$ws.Cells.Item(16,10).Value = '{"Country": "Deutschland", "City": "Gerbrunn", "ZipCode": "97218", "Street": "Helene Wessel Weg 5", "LocationId": 12047, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(16,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "3 Kaltmieten", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in Gerbrunn", "Link": "https://www.immowelt.de/immobilienpreise/gerbrunn/mietspiegel"}}, "DataTable": [{"NumberValue": 2880, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 400, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 3280, "Unit": "EUR"}]}'
$ws.Cells.Item(16,12).Value = '["https://ms.immowelt.org/17e5aca9-793d-40b3-aab5-a5dcea349a87/0d9236ef-f790-4bac-82dc-9a229cfe4527", "https://ms.immowelt.org/cf96b674-5a06-4156-b363-d3e23b5dfe92/c12ffbff-3538-4b01-a604-ee0869a952b6", "https://ms.immowelt.org/2755511b-7b13-42f8-bad0-3d769da9305b/b3d17421-aa79-4be1-83cb-ae6fe7955297", "https://ms.immowelt.org/17109eb9-b8cf-435e-b0b9-8377b8a7c9f4/986295b9-a109-4924-b6b8-fc375860cd30", "https://ms.immowelt.org/d4384af4-fbdc-4807-9ab2-e4c02ad5e5b9/61dd1708-3529-4a40-9890-401b9917edfc", "https://ms.immowelt.org/efa4e805-be82-4e38-998b-46485660f4e5/fbb8ea8f-2b17-45b1-8691-eb2fe09dd0a7"]'
$ws.Cells.Item(16,13).Value = '""'
$ws.Cells.Item(16,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Traumwohnung in Top Lage mit toller Weitsicht!", "basicObjectPricEur": 2880, "basicLivingSpace": 195, "basicRooms": 6, "basicConstructionYear": 2023, "basicCity": "Gerbrunn", "basicStreet": "Helene Wessel Weg 5", "basicUrl": "https://www.immowelt.de/expose/2agbf5s", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 17
$ws.Cells.Item(17,1).Value = '"https://www.immowelt.de/expose/2a2yn5v"'
$ws.Cells.Item(17,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(17,3).Value = '"Optimal f\u00fcr Studenten-WG : 2-Zimmerwohnung in zentraler Lage im Frauenland"'
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '580'
$ws.Cells.Item(17,5).Value = '"EUR"'
$ws.Cells.Item(17,6).NumberFormat = "@"
$ws.Cells.Item(17,6).Value = '44'
$ws.Cells.Item(17,7).NumberFormat = "@"
$ws.Cells.Item(17,7).Value = '2'
$ws.Cells.Item(17,8).Value = 'null'
$ws.Cells.Item(17,9).Value = '["WANNE"]'
$ws.Cells.Item(17,10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Altstadt", "ZipCode": "97072", "LocationId": 496004, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(17,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1.080,00 ", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Altstadt)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-altstadt/mietspiegel"}}, "DataTable": [{"NumberValue": 580, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 120, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Warmmiete enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}]}'
$ws.Cells.Item(17,12).Value = '["https://ms.immowelt.org/67968400-d13b-449c-8f0f-1f8dc3645f4c/4d6bfee6-a4bb-438e-a524-96840bd51ae1", "https://ms.immowelt.org/c3319402-5904-4d76-95d2-db5b70bee409/4e68c2f5-c18e-4b63-bc8f-0ee1c1309ebf", "https://ms.immowelt.org/2dd83638-fb99-4566-ba6f-cf7b62d7c4bb/6ceb0ef0-07a4-468f-8f2a-f2e6741a577f", "https://ms.immowelt.org/1609ec66-fe8f-4f74-b4d5-02421f7513b0/a1757c96-7393-44bc-bf2e-6c9466ce4f45", "https://ms.immowelt.org/35e2e29e-1627-4973-99d4-4c3b1facdc4e/b477018c-9c03-4c1b-861c-11a14bcff4e8", "https://ms.immowelt.org/48e33aa4-fbdd-48d1-ac41-d25502baa7d5/14df5c1a-b256-4cd0-a73b-72d4dd9a2961"]'
$ws.Cells.Item(17,13).Value = '{"companyName": "Viefhaus Immobilien Inh. Sebastian Viefhaus", "address": {"city": "W\u00fcrzburg", "zipCode": "97074", "street": "Fichtestra\u00dfe 10"}, "salutation": "Herr", "firstName": "Thomas", "lastName": "Meister", "phone": "0931-88065061"}'
$ws.Cells.Item(17,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Optimal f\u00fcr Studenten-WG : 2-Zimmerwohnung in zentraler Lage im Frauenland", "basicObjectPricEur": 580, "basicLivingSpace": 44, "basicRooms": 2, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2a2yn5v", "basicContactPhone": "0931-88065061"}'

# Row 18
$ws.Cells.Item(18,1).Value = '"https://www.immowelt.de/expose/2upu84s"'
$ws.Cells.Item(18,2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(18,3).Value = '"Gro\u00dfe 3-Zimmer-Wohnung mit eigenem Garten und Balkon in K\u00fcrnach"'
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '900'
$ws.Cells.Item(18,5).Value = '"EUR"'
$ws.Cells.Item(18,6).NumberFormat = "@"
$ws.Cells.Item(18,6).Value = '90'
$ws.Cells.Item(18,7).NumberFormat = "@"
$ws.Cells.Item(18,7).Value = '3'
$ws.Cells.Item(18,8).Value = 'null'
$ws.Cells.Item(18,9).Value = '["WANNE", "FENSTER", "Kelleranteil", "ELEKTRO", "Garten", "gartennutzung", "Balkon", "Terrasse"]'
$ws.Cells.Item(18,10).Value = '{"Country": "Deutschland", "City": "K\u00fcrnach", "ZipCode": "97273", "LocationId": 12666, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(18,11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "2 Kaltmieten ", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in K\u00fcrnach", "Link": "https://www.immowelt.de/immobilienpreise/kuernach/mietspiegel"}}, "DataTable": [{"NumberValue": 900, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 80, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 50, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "1 Stellplatz"}]}'
$ws.Cells.Item(18,12).Value = '["https://ms.immowelt.org/11008e51-5945-4ab3-a0eb-f2a703785b20/53fa6be7-ade6-4395-81ca-d66b9bbb1ba9", "https://ms.immowelt.org/64bc93db-da55-44ca-a2c8-c14f8ce86353/14ba8183-9d5b-4c81-a957-284d595a032f", "https://ms.immowelt.org/87e6208f-ee9e-4d96-a62a-6561d184b0e3/572c15dc-ed2e-4296-bc44-96a28e7cf135", "https://ms.immowelt.org/ba824adc-5910-4328-b432-ac6c69a928b3/dc0e3963-adad-4942-ad9a-dc0fc1b71796", "https://ms.immowelt.org/cee10f88-e7a3-499e-be76-626c1d42e77d/c1805001-8c80-4257-977c-80d06c2cbbc7", "https://ms.immowelt.org/aba27baa-44fb-49d2-833b-ffe9285cf719/2f9875be-ad3f-435e-bf3a-d7d63a5e8a17", "https://ms.immowelt.org/b63c1740-9745-464b-8364-2df590a1d676/f71392e8-2cac-4fad-aeaf-0379e7187348", "https://ms.immowelt.org/0b219ac9-8c45-4c18-9228-d206e9ff703e/4bce456f-aecf-48fa-a4a2-b0b332e7badc", "https://ms.immowelt.org/711ca5fb-12d7-42dc-a1f6-32893287c0eb/4002136b-ce99-43ba-8406-dee393492ff0", "https://ms.immowelt.org/17453203-a015-4bcc-816b-12fe9fe761ee/3c33ba0b-8a93-4d03-b3c7-0267e690bc76", "https://ms.immowelt.org/f09d3348-84eb-4a67-941e-68eee061b36d/75a81b43-153a-44e2-8336-9f323ccdad5a", "https://ms.immowelt.org/79e28d05-fb4a-4b73-ad94-cd78318ef65b/0441105e-ad78-4d19-86b1-dab1c3c87ee1", "https://ms.immowelt.org/07c2164c-067e-490d-b769-760ef64271e2/2e82cebf-00da-4337-83e8-6186da964d0d", "https://ms.immowelt.org/de5a7379-e440-43c4-8670-4bdd1e8bdca0/e80fa6cc-bbb1-48ab-b091-2b052e63948b", "https://ms.immowelt.org/9cc41d7e-f16c-4d38-ae23-e1f15da29cfc/e1740437-8eb5-4984-bba2-1c1036d06b98"]'
$ws.Cells.Item(18,13).Value = '""'
$ws.Cells.Item(18,14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Gro\u00dfe 3-Zimmer-Wohnung mit eigenem Garten und Balkon in K\u00fcrnach", "basicObjectPricEur": 900, "basicLivingSpace": 90, "basicRooms": 3, "basicConstructionYear": null, "basicCity": "K\u00fcrnach", "basicUrl": "https://www.immowelt.de/expose/2upu84s", "basicContactPhone": "", "basicContactMobile": ""}'
